# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AC, AD, AE
$ws.Range("AC1").Value2 = "Wins"
$ws.Range("AD1").Value2 = "Losses"
$ws.Range("AE1").Value2 = "Ties"

# Match header formatting/style used by the other header cells (e.g. A1)
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-48: Wins = 61, Losses = 101, Ties = 0 for every row
$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value2 = 61   # AC
    $ws.Cells.Item($r, 30).Value2 = 101  # AD
    $ws.Cells.Item($r, 31).Value2 = 0    # AE
}
